$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Студент"
$ws.Range("B1").Value = "Тема проеткту"
$ws.Range("A2").Value = "Вохранов І.А."
$ws.Range("B2").Value = "Система комунікації орендодавця та орендатора"

$ws.Columns.Item(1).ColumnWidth = 44.44
$ws.Columns.Item(2).ColumnWidth = 71.44
